$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to nombre_aides (column C) and montant_total (column E)
# for the 2022-06-24 data refresh of the Fonds de solidarite dataset.
$updates = @(
    @{ Row = 3;   C = 249336;  E = 1036489581 },
    @{ Row = 53;  C = 141688;  E = 590078732 },
    @{ Row = 91;  C = 151249;  E = 483346675 },
    @{ Row = 92;  C = 409341;  E = 1597795335 },
    @{ Row = 93;  C = 209683;  E = 1310346798 },
    @{ Row = 94;  C = 94246;   E = 919357089 },
    @{ Row = 95;  C = 50813;   E = 935119531 },
    @{ Row = 96;  C = 17331;   E = 797952512 },
    @{ Row = 104; C = 135345;  E = 273027754 },
    @{ Row = 116; C = 4567;    E = 20680072 },
    @{ Row = 118; C = 981;     E = 11896140 },
    @{ Row = 145; C = 11834;   E = 182736103 },
    @{ Row = 163; C = 70985;   E = 131769402 },
    @{ Row = 173; C = 96866;   E = 327978644 },
    @{ Row = 174; C = 226115;  E = 900825976 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
